$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.694.94"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "1.893.24"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").Value = "'311.48"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").Value = "'0.4948"
$ws.Range("E7").Value = "  +2.06%  "
$ws.Range("D8").Value = "'0.3794"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("D9").Value = "'0.07323"
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("D10").Value = "'0.9085"
$ws.Range("E10").Value = "  -4.75%  "
$ws.Range("D11").Value = "'20.61"
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("D12").Value = "'0.07643"
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("D13").Value = "1.865.98"
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("D14").Value = "'5.469"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").Value = "'6.640"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").Value = "'91.12"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "'0.000008737"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").Value = "27.651.57"
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("D21").Value = "'14.46"
$ws.Range("E21").Value = "  -4.18%  "
$ws.Range("D22").Value = "'5.119"
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").Value = "2.112.10"
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("D25").Value = "'154.24"
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("D26").Value = "'1.840"
$ws.Range("E26").Value = "  -5.69%  "
$ws.Range("D27").Value = "'18.39"
$ws.Range("E27").Value = "  -1.31%  "
$ws.Range("D28").Value = "'2.169"
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("D29").Value = "'115.34"
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("D30").Value = "'4.872"
$ws.Range("E30").Value = "  -3.61%  "
$ws.Range("D31").Value = "'0.08941"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "'3.202"
$ws.Range("E32").Value = "  -4.36%  "
$ws.Range("D33").Value = "'1.225"
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("D34").Value = "'0.7669"
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("D35").Value = "'4.633"
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("D36").Value = "'2.562"
$ws.Range("E36").Value = "  -8.00%  "
$ws.Range("D37").Value = "'0.02041"
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("D38").Value = "'1.096"
$ws.Range("E38").Value = "  -3.36%  "
$ws.Range("D39").Value = "'0.05282"
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("D40").Value = "'0.5488"
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("D41").Value = "'2.987"
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("D42").Value = "'6.883"
$ws.Range("E42").Value = "  -3.73%  "
$ws.Range("D43").Value = "'8.540"
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").Value = "'112.60"
$ws.Range("E44").Value = "  +4.46%  "
$ws.Range("D45").Value = "'0.1519"
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("D46").Value = "'10.62"
$ws.Range("E46").Value = "  -2.02%  "
$ws.Range("D47").Value = "'0.4798"
$ws.Range("E47").Value = "  -3.53%  "
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("D49").Value = "'1.635"
$ws.Range("E49").Value = "  -3.01%  "
$ws.Range("D50").Value = "'67.31"
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("E51").Value = "  -1.46%  "

Write-Output "Applied all changes"
